$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2 ("Update: Plan for the project") - content placeholder (shape 2)
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(2)
$tf2 = $shape2.TextFrame
$tr2 = $tf2.TextRange

# bodyPr: normAutofit -> noAutofit
$tf2.AutoSize = 0

# Insert a brand-new first paragraph ("Goal: ...") ahead of "Topic: Ray Tracing "
$tr2.InsertBefore("Goal:  designed for showing the strength of ray tracing and customize a set of attributes in Cycles Render Engine particular for food rendering.`r")

# --- paragraph 1 ("Goal: ...") : split into its distinctly-formatted runs ---
# "Goal"
$r = $tr2.Characters(1,4)
$r.Font.Size = 20
$r.Font.Bold = 1
$r.Font.Name = "+mj-lt"
# ": "
$r = $tr2.Characters(5,2)
$r.Font.Size = 20
$r.Font.Name = "+mj-lt"
# " designed for showing the strength of ray tracing and customize a set of attributes in "
$r = $tr2.Characters(7,87)
$r.Font.Size = 20
$r.Font.Bold = 0
$r.Font.Italic = 0
$r.Font.Shadow = 1
$r.Font.Shadow = 0
$r.Font.Name = "+mj-lt"
# "Cycles "
$r = $tr2.Characters(94,7)
$r.Font.Size = 20
$r.Font.Name = "+mj-lt"
# "Render Engine "
$r = $tr2.Characters(101,14)
$r.Font.Size = 20
$r.Font.Name = "+mj-lt"
# "particular "
$r = $tr2.Characters(115,11)
$r.Font.Size = 20
$r.Font.Bold = 0
$r.Font.Italic = 0
$r.Font.Shadow = 1
$r.Font.Shadow = 0
$r.Font.Name = "+mj-lt"
# "for food rendering."
$r = $tr2.Characters(126,19)
$r.Font.Size = 20
$r.Font.Bold = 0
$r.Font.Italic = 0
$r.Font.Shadow = 1
$r.Font.Shadow = 0
$r.Font.Name = "+mj-lt"

# --- paragraph 2 ("Topic: Ray Tracing ") : level 0 -> gets sz=2000 + latin ---
$para = $tr2.Paragraphs(2)
$para.Font.Size = 20
$para.Font.Name = "+mj-lt"

# --- paragraph 3 ("Still related...") : level 1, no size change, just latin ---
$para = $tr2.Paragraphs(3)
$para.Font.Name = "+mj-lt"

# --- paragraph 4 ("Cycles" + " in Blender...") : level 1, latin on both runs ---
$r = $tr2.Characters($tr2.Paragraphs(4).Start, 6)
$r.Font.Name = "+mj-lt"
$r = $tr2.Characters($tr2.Paragraphs(4).Start + 6, $tr2.Paragraphs(4).Length - 6)
$r.Font.Name = "+mj-lt"

# --- paragraph 5 ("Understand Ray Tracing...") : level 1, just latin ---
$para = $tr2.Paragraphs(5)
$para.Font.Name = "+mj-lt"

# --- paragraph 6 ("Customize Ray Tracing...") : level 1, just latin ---
$para = $tr2.Paragraphs(6)
$para.Font.Name = "+mj-lt"

# --- paragraph 7 ("Build a scene...") : level 0 -> sz=2000 + latin ---
$para = $tr2.Paragraphs(7)
$para.Font.Size = 20
$para.Font.Name = "+mj-lt"

# --- paragraph 8 ("San Miguel Scene...") : level 1, just latin ---
$para = $tr2.Paragraphs(8)
$para.Font.Name = "+mj-lt"

# --- paragraph 9 ("Too complicated...") : level 2 -> sz=2000 + latin ---
$para = $tr2.Paragraphs(9)
$para.Font.Size = 20
$para.Font.Name = "+mj-lt"

# --- paragraph 10 ("Besides donuts...") : level 1, just latin ---
$para = $tr2.Paragraphs(10)
$para.Font.Name = "+mj-lt"

# --- paragraph 11 ("Plan to modeling...") : level 2 -> sz=2000 + latin ---
$para = $tr2.Paragraphs(11)
$para.Font.Size = 20
$para.Font.Name = "+mj-lt"

# --- paragraph 13 (final, empty endParaRPr lang="zh-CN") : sz=2000 only ---
$para = $tr2.Paragraphs(13)
$para.Font.Size = 20

# ---------------------------------------------------------------------------
# Slide 5 ("Takeaways for this week") - content placeholder (shape 2)
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shape5 = $s5.Shapes.Item(2)
$tr5 = $shape5.TextFrame.TextRange

$tr5.Paragraphs(1).Text = "Blender Add-on"
$tr5.Paragraphs(2).Text = "Blender directory structures"
$tr5.InsertAfter("`rpython scripting")
$tr5.Paragraphs(3).Font.Size = 36
